# Update cryptos list cell values (text) to match latest scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '68.368.22'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +0.07%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.646.10'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +0.21%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '598.57'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +0.24%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '158.77'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +2.80%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.541'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -0.87%  '
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -0.80%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.156'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -1.07%  '
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +0.64%  '
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -0.74%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '27.96'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -0.49%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '3.131.97'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +0.39%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.0000187'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -2.47%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '68.286.95'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +0.08%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.604.32'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -1.54%  '
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -0.16%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '359.69'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -1.03%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.40'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -1.16%  '
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +0.97%  '
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -2.43%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '2.07'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +0.43%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '74.38'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -0.36%  '
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +0.07%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.79'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -0.54%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.779.82'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -0.05%  '
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -2.54%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.999'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -0.26%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '562.37'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -1.72%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '8.03'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -1.38%  '
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -2.83%  '
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +1.27%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.65'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +3.59%  '
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +0.03%  '
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -1.45%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '160.26'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -0.49%  '
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +1.55%  '
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -0.83%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.87'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -1.49%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.32'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -1.16%  '
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -1.36%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0₆0322'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -3.74%  '
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +0.02%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '157.37'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +0.67%  '
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +0.72%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '22.00'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +0.66%  '
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -1.42%  '
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -2.05%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.574'
